$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Move the active selection from F8 to B11.
$ws.Range("B11").Select() | Out-Null

# E9 keeps its value (19 May 2017 / 42874) but picks up the date number
# format already used by E16 (style index 4), instead of the unformatted
# "General" style it had before. Copy/PasteSpecial(formats) transfers the
# format without touching the underlying value.
$ws.Range("E16").Copy()
$ws.Range("E9").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Push the deadline dates for the "wireframe" block (rows 15, 16, 18, 19)
# from 19 May 2017 forward one week to 26 May 2017.
$ws.Range("E15").Value = 42881
$ws.Range("E16").Value = 42881
$ws.Range("E18").Value = 42881
$ws.Range("E19").Value = 42881

# Push the deadline dates for the "mockup/realisatie" block (rows 21-25)
# from 26 May 2017 forward one week to 2 June 2017.
$ws.Range("E21").Value = 42888
$ws.Range("E22").Value = 42888
$ws.Range("E23").Value = 42888
$ws.Range("E24").Value = 42888
$ws.Range("E25").Value = 42888
